# Append: 2026-02-06 06:52 JST
# Refresh the scraped Lancers listing on the "ランサーズ" sheet with a new
# batch of results (12 data rows instead of 13), new timestamps, and
# slightly narrower/wider columns B, D, H.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Drop the two rows that fell off the bottom of the new batch (13 & 14) ---
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(13).Delete()

# --- 2. Clear every existing hyperlink; we'll re-add exactly the ones needed ---
$ws.Cells.Hyperlinks.Delete()

# --- 3. Overwrite rows 2-12 with the new scrape results ---

# Row 2
$ws.Range("A2").Value = "2026-02-06 06:52:00"
$ws.Range("B2").Value = "【急募】Next.js × Expoでアプリ開発仲間を大募集!"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5486863"
$ws.Range("G2").Value = 218
$ws.Range("H2").Value = "🔥Next.js ◆開発 ◇アプリ"

# Row 3
$ws.Range("A3").Value = "2026-02-06 06:52:00"
$ws.Range("B3").Value = "【急募】WordPressサイト再構築+LINE・予約連携+顧客管理機能構築|テーマ指定あり|"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5486856"
$ws.Range("G3").Value = 93
$ws.Range("H3").Value = "◇サイト ○WordPress"

# Row 4
$ws.Range("A4").Value = "2026-02-06 06:52:00"
$ws.Range("B4").Value = "【業務改善】訪問業務に特化したスケジュール/介入実績管理Webシステム構築"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5486583"
$ws.Range("G4").Value = 85
$ws.Range("H4").Value = "◇業務改善"

# Row 5
$ws.Range("A5").Value = "2026-02-06 06:52:00"
$ws.Range("B5").Value = "【Java/講師】新入社員研修のサブ講師募集"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5477144"
$ws.Range("G5").Value = 78
$ws.Range("H5").Value = "★Java"

# Row 6
$ws.Range("A6").Value = "2026-02-06 06:52:00"
$ws.Range("B6").Value = "初心者向けダンススクールの問い合わせフォームを置き換える/拡張するチャットボット開発"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5486851"
$ws.Range("G6").Value = 75
$ws.Range("H6").Value = "◆開発"

# Row 7
$ws.Range("A7").Value = "2026-02-06 06:52:00"
$ws.Range("B7").Value = "【募集】PHP + MySQLでのcron用スクリプト作成依頼"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5487013"
$ws.Range("G7").Value = 50
$ws.Range("H7").Value = "◇MySQL ○PHP"

# Row 8 (no skill-summary column this time)
$ws.Range("A8").Value = "2026-02-06 06:52:00"
$ws.Range("B8").Value = "【長期】寝具ブランドのAmazon・楽天市場 運用代行パートナー募集"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5486471"
$ws.Range("G8").Value = 25
$ws.Range("H8").ClearContents()

# Row 9
$ws.Range("A9").Value = "2026-02-06 06:52:00"
$ws.Range("B9").Value = "【急募】外部CTOを探しています!"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5486956"
$ws.Range("G9").Value = 18
$ws.Range("H9").ClearContents()

# Row 10
$ws.Range("A10").Value = "2026-02-06 06:52:00"
$ws.Range("B10").Value = "【急募】SSLエラー解決のための専門家を探しています"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5486960"
$ws.Range("G10").Value = 13
$ws.Range("H10").ClearContents()

# Row 11
$ws.Range("A11").Value = "2026-02-06 06:52:00"
$ws.Range("B11").Value = "【急募】Klaviyoスパム対策とドメイン解決の専門家募集"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5486673"
$ws.Range("G11").Value = 13
$ws.Range("H11").ClearContents()

# Row 12
$ws.Range("A12").Value = "2026-02-06 06:52:00"
$ws.Range("B12").Value = "【急募】BOXファイルをGASでkintoneに自動同期したい"
$ws.Range("C12").Value = "システム開発"
$ws.Range("D12").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E12").Value = "期限情報なし"
$ws.Range("F12").Value = "https://www.lancers.jp/work/detail/5487010"
$ws.Range("G12").Value = 10
$ws.Range("H12").ClearContents()

# --- 4. Re-create hyperlinks for F2:F12, pointing at the new URLs ---
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5486863")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5486856")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5486583")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5477144")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5486851")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5487013")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5486471")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5486956")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5486960")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5486673")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5487010")

# --- 5. Column-width tweaks (B 51->49, D 32->30, H 18->19) ---
# ColumnWidth is in "characters"; the engine adds a constant ~0.8333 padding
# on top of whatever is set, so back that out to land on the exact target.
$pad = 0.8333333333333334
$ws.Columns.Item(2).ColumnWidth = 49 - $pad
$ws.Columns.Item(4).ColumnWidth = 30 - $pad
$ws.Columns.Item(8).ColumnWidth = 19 - $pad
